# Add a new "Logout redirect" entry to the simpleXWiki modifications list on
# Sheet1, right above the existing "menuview.vm" / "Menu View Extras" row.
# This mirrors the source commit:
#   "menuview.vm modified to include logoutRedirect variable for controlling
#    default login redirect page"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The existing "menuview.vm" entry lives on row 34. Insert a brand new row
# above it so the new "Logout redirect" variable becomes the first entry
# listed under the menuview.vm template, pushing the old "Menu View Extras"
# row (and everything below it) down by one.
$ws.Rows.Item(34).Insert()

$newRow = 34

$ws.Cells.Item($newRow, 1).Value = "menuview.vm"
$ws.Cells.Item($newRow, 2).Value = "Logout redirect"
$ws.Cells.Item($newRow, 3).Value = "Option to include additional content to the left of the rightmenu content (profile menu)."
$ws.Cells.Item($newRow, 4).Value = "logoutRedirect"
$ws.Cells.Item($newRow, 5).Value = "Variable replacement"
$ws.Cells.Item($newRow, 6).Value = "Replace the value of the xredirect paramter from `$xwiki.relativeRequestURL to the value of `$logoutRedirect"

# The "menuview.vm" label in column A only appears once per template group;
# clear it from what is now row 35 (the original "Menu View Extras" row)
# so the grouping matches the rest of the sheet.
$ws.Cells.Item(35, 1).Value = ""

# Update the selection to roughly match where the author was last working.
$ws.Range("G34").Select()
